$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row: "<name>_old" -> "<name>_FV2210" (cols A-J)
#        and "<name>_new" -> "<name>_FV2304" (cols L-U). Column K ("diff") unchanged.
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $baseNames[$i] + "_FV2304"
}

# --- 2) Freeze the header row (ySplit=1, top-left cell A2)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Turn the used range into an Excel Table ("Table1") with autofilter + banded rows
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U77"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false
